$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.696.02"
$ws.Range("E2").Value = "  -3.62%  "
$ws.Range("D3").Value = "'2.280.84"
$ws.Range("E3").Value = "  -4.33%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'530.51"
$ws.Range("E5").Value = "  -4.39%  "
$ws.Range("D6").Value = "'129.62"
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.578"
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("D9").Value = "'2.279.40"
$ws.Range("E9").Value = "  -4.28%  "
$ws.Range("D10").Value = "'0.0992"
$ws.Range("E10").Value = "  -5.24%  "
$ws.Range("D11").Value = "'5.39"
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'0.328"
$ws.Range("E13").Value = "  -4.52%  "
$ws.Range("D14").Value = "'23.34"
$ws.Range("E14").Value = "  -4.38%  "
$ws.Range("D15").Value = "'2.684.75"
$ws.Range("E15").Value = "  -4.41%  "
$ws.Range("D16").Value = "'57.626.00"
$ws.Range("D17").Value = "'0.0000131"
$ws.Range("E17").Value = "  -4.25%  "
$ws.Range("D18").Value = "'2.275.98"
$ws.Range("E18").Value = "  -4.45%  "
$ws.Range("D19").Value = "'10.45"
$ws.Range("E19").Value = "  -5.88%  "
$ws.Range("D20").Value = "'4.20"
$ws.Range("E20").Value = "  -6.42%  "
$ws.Range("D21").Value = "'312.37"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").Value = "'6.34"
$ws.Range("E22").Value = "  -5.54%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'62.46"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("E25").Value = "  -4.30%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'7.97"
$ws.Range("E27").Value = "  -5.32%  "
$ws.Range("E28").Value = "  -6.23%  "
$ws.Range("D29").Value = "'170.18"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("E30").Value = "  -5.88%  "
$ws.Range("D31").Value = "'0.0₃0714"
$ws.Range("E31").Value = "  -5.86%  "
$ws.Range("E32").Value = "  -5.47%  "
$ws.Range("D33").Value = "'1.03"
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("D34").Value = "'0.376"
$ws.Range("E34").Value = "  -5.30%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'17.70"
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "'1.23"
$ws.Range("E38").Value = "  -6.69%  "
$ws.Range("E39").Value = "  -6.52%  "
$ws.Range("D40").Value = "'37.95"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("E41").Value = "  -5.86%  "
$ws.Range("D42").Value = "'140.42"
$ws.Range("E42").Value = "  -3.64%  "
$ws.Range("D43").Value = "'286.35"
$ws.Range("E43").Value = "  -10.01%  "
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("D45").Value = "'0.0945"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("E46").Value = "  -3.16%  "
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("D48").Value = "'17.99"
$ws.Range("E48").Value = "  -8.60%  "
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").Value = "'0.0₆0200"
$ws.Range("E51").Value = "  +84.25%  "
